$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 10 corresponds to
# cd2cbb82-99d4-4be8-83d2-a81fa423dc03.96fb3de54468e89b7de6783c033bbe4bef4ff415.zh-cn.xlf
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D10").Value = "2016-03-01 09:20:36"
$wsZh.Range("G10").Value = "2016-03-01 09:21:19"

# de-de sheet: row 10 corresponds to
# cd2cbb82-99d4-4be8-83d2-a81fa423dc03.96fb3de54468e89b7de6783c033bbe4bef4ff415.de-de.xlf
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D10").Value = "2016-03-01 09:20:46"
$wsDe.Range("G10").Value = "2016-03-01 09:21:37"
